# Auto-generated script applying cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. '0.138', '62.119.92') are written back as text, matching the
# original inline-string cell type instead of being coerced to numbers.
$dataRange = $ws.Range('B2:E51')
$dataRange.NumberFormat = '@'

$ws.Range('D2').Value = '62.119.92'
$ws.Range('E2').Value = '  -6.97%  '
$ws.Range('D3').Value = '2.922.84'
$ws.Range('E3').Value = '  -9.04%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '530.12'
$ws.Range('E5').Value = '  -10.89%  '
$ws.Range('D6').Value = '128.31'
$ws.Range('E6').Value = '  -15.52%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '2.899.30'
$ws.Range('E8').Value = '  -9.47%  '
$ws.Range('E9').Value = '  -17.81%  '
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  -19.59%  '
$ws.Range('D11').Value = '5.73'
$ws.Range('E11').Value = '  -11.48%  '
$ws.Range('D12').Value = '0.419'
$ws.Range('E12').Value = '  -15.35%  '
$ws.Range('D13').Value = '30.99'
$ws.Range('E13').Value = '  -20.56%  '
$ws.Range('D14').Value = '0.0000196'
$ws.Range('E14').Value = '  -19.72%  '
$ws.Range('D15').Value = '3.399.06'
$ws.Range('E15').Value = '  -9.03%  '
$ws.Range('D16').Value = '62.093.82'
$ws.Range('E16').Value = '  -7.03%  '
$ws.Range('E17').Value = '  -4.89%  '
$ws.Range('D18').Value = '2.930.56'
$ws.Range('E18').Value = '  -8.90%  '
$ws.Range('D19').Value = '460.53'
$ws.Range('E19').Value = '  -13.18%  '
$ws.Range('D20').Value = '6.01'
$ws.Range('E20').Value = '  -15.70%  '
$ws.Range('D21').Value = '12.39'
$ws.Range('E21').Value = '  -16.65%  '
$ws.Range('D22').Value = '0.617'
$ws.Range('E22').Value = '  -18.72%  '
$ws.Range('D23').Value = '6.24'
$ws.Range('E23').Value = '  -21.41%  '
$ws.Range('D24').Value = '73.54'
$ws.Range('E24').Value = '  -14.24%  '
$ws.Range('D25').Value = '11.53'
$ws.Range('E25').Value = '  -17.19%  '
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = '2.59'
$ws.Range('E27').Value = '  -18.78%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '1.79'
$ws.Range('E28').Value = '  -17.78%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '6.69'
$ws.Range('E29').Value = '  -17.75%  '
$ws.Range('D30').Value = '24.04'
$ws.Range('E30').Value = '  -17.47%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('B32').Value = 'Mantle'
$ws.Range('C32').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D32').Value = '1.04'
$ws.Range('E32').Value = '  -9.96%  '
$ws.Range('D33').Value = '2.27'
$ws.Range('E33').Value = '  -15.73%  '
$ws.Range('D34').Value = '50.70'
$ws.Range('E34').Value = '  -5.47%  '
$ws.Range('D35').Value = '451.45'
$ws.Range('E35').Value = '  -17.12%  '
$ws.Range('D36').Value = '5.29'
$ws.Range('E36').Value = '  -18.29%  '
$ws.Range('D37').Value = '4.51'
$ws.Range('E37').Value = '  -21.15%  '
$ws.Range('D38').Value = '0.0376'
$ws.Range('E38').Value = '  -11.94%  '
$ws.Range('D39').Value = '0.0731'
$ws.Range('E39').Value = '  -14.79%  '
$ws.Range('E40').Value = '  -12.47%  '
$ws.Range('D41').Value = '7.54'
$ws.Range('E41').Value = '  -18.95%  '
$ws.Range('D42').Value = '2.583.90'
$ws.Range('E42').Value = '  -11.03%  '
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = '2.12'
$ws.Range('E44').Value = '  -20.84%  '
$ws.Range('D45').Value = '0.212'
$ws.Range('E45').Value = '  -19.79%  '
$ws.Range('D46').Value = '108.60'
$ws.Range('E46').Value = '  -7.94%  '
$ws.Range('D47').Value = '0.0979'
$ws.Range('E47').Value = '  -14.38%  '
$ws.Range('D48').Value = '1.72'
$ws.Range('E48').Value = '  -20.32%  '
$ws.Range('D49').Value = '0.0₃0445'
$ws.Range('E49').Value = '  -23.43%  '
$ws.Range('B50').Value = 'BitgetToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D50').Value = '1.15'
$ws.Range('E50').Value = '  -7.55%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '20.62'
$ws.Range('E51').Value = '  -22.39%  '

# Remove the temporary text formatting so styling matches the original
# (unstyled data cells with no explicit 's' attribute).
$dataRange.ClearFormats()
